$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.355.24'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '3.687.86'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '679.54'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '159.17'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("E9").Value = '  -1.06%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.36%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Value = '4.308.94'
$ws.Range("E13").Value = '  -0.06%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '32.41'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("D15").Value = '3.693.46'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '69.312.16'
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("E17").Value = '  +2.78%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '15.98'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E19").Value = '  -0.82%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '467.80'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("E22").Value = '  -0.81%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '80.00'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '3.833.52'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -0.07%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0000123'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -5.08%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.90'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -1.99%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -3.84%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.57'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.35%  '
$ws.Range("E32").Value = '  -3.10%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").Value = '3.676.83'
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("E36").Value = '  -4.69%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '8.34'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.42%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  -0.06%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.14%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0903'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.95%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '170.76'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +4.33%  '
$ws.Range("E44").Value = '  -1.24%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '47.41'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.77%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '28.18'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -5.44%  '
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.30'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.14%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.000276'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("E51").Value = '  -3.01%  '
